$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated "mL Filt" (column C) values -> dependent formulas in D, J, K, M, N
# recalc automatically since they reference column C (via column D).
$ws.Range("C55").Value = 160
$ws.Range("C64").Value = 120
$ws.Range("C71").Value = 120
$ws.Range("C91").Value = 150
$ws.Range("C101").Value = 150

# Freeze the header rows (rows 1-7) and select column B, matching the
# reviewer's view state when re-opening the notebook.
$ws.Range("A8").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B:B").Select()
